# Ratios.xlsx - "corrected age calculation algorithm" fix:
#  - dU234 (col A) and its absolute error (col B) are recomputed with the
#    corrected algorithm (mean/median switch for outlier detection) for
#    every data row (rows 2-16 on the "Ratios" sheet).
#  - column A got one character-unit narrower (21.71 -> 20.71 OOXML width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A by ~1 character. The host quantises COM ColumnWidth to
# whole display pixels, so 19.8333... (-> stored OOXML width 20.666667) is
# the closest representable value to the target 20.7109375.
$ws.Columns.Item(1).ColumnWidth = 19.833333333333332

$ws.Range("A2").Value = -3.543703889504424
$ws.Range("B2").Value = 0.001295645475698015

$ws.Range("A3").Value = -0.4998166906194923
$ws.Range("B3").Value = 0.001462523696273757

$ws.Range("A4").Value = -5.0778265256608
$ws.Range("B4").Value = 0.001112384267381515

$ws.Range("A5").Value = -0.0622149226449098
$ws.Range("B5").Value = 0.001651962676204683

$ws.Range("A6").Value = -2.259540569212626
$ws.Range("B6").Value = 0.0009824633009997295

$ws.Range("A7").Value = -2.29774101750313
$ws.Range("B7").Value = 0.0007393942454393191

$ws.Range("A8").Value = -3.36255752403547
$ws.Range("B8").Value = 0.0008344977463790538

$ws.Range("A9").Value = -1.280644752506777
$ws.Range("B9").Value = 0.0006234212432108183

$ws.Range("A10").Value = -4.023575486949205
$ws.Range("B10").Value = 0.001047608192131502

$ws.Range("A11").Value = -1.400267488405271
$ws.Range("B11").Value = 0.0006150742400669406

$ws.Range("A12").Value = -4.00234612944772
$ws.Range("B12").Value = 0.0008888145216077643

$ws.Range("A13").Value = -2.14492949082945
$ws.Range("B13").Value = 0.000679718985948966

$ws.Range("A14").Value = -3.699781135822144
$ws.Range("B14").Value = 0.0007902626358170166

$ws.Range("A15").Value = -1.642549567724694
$ws.Range("B15").Value = 0.0005192005057607072

$ws.Range("A16").Value = -4.958884621684923
$ws.Range("B16").Value = 0.001016860524648371
